$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a handful of D-column cells to remain Text (not auto-converted to Number),
# so trailing zeros in values like 0.0960 / 1.20 are preserved exactly, matching the source data.
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

# Apply the updated cryptos snapshot values.
$ws.Range("D2").Value = '43.412.02'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '2.344.99'
$ws.Range("E3").Value = '  +3.02%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '232.53'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '0.649'
$ws.Range("E6").Value = '  +2.01%  '
$ws.Range("D7").Value = '65.83'
$ws.Range("E7").Value = '  +3.55%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.453'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").Value = '0.0960'
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("D11").Value = '56.92'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").Value = '26.79'
$ws.Range("E12").Value = '  -2.02%  '
$ws.Range("D13").Value = '2.693.40'
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '15.41'
$ws.Range("E15").Value = '  -2.29%  '
$ws.Range("D16").Value = '6.26'
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").Value = '0.848'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").Value = '2.343.71'
$ws.Range("E18").Value = '  +3.49%  '
$ws.Range("D19").Value = '43.366.40'
$ws.Range("E19").Value = '  -1.22%  '
$ws.Range("D20").Value = '0.0₃0981'
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("D21").Value = '74.16'
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").Value = '6.22'
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("D23").Value = '249.27'
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("D24").Value = '3.85'
$ws.Range("E24").Value = '  +16.47%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("D28").Value = '9.91'
$ws.Range("E28").Value = '  -1.38%  '
$ws.Range("D29").Value = '174.86'
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("D30").Value = '22.17'
$ws.Range("E30").Value = '  +6.10%  '
$ws.Range("E31").Value = '  +6.58%  '
$ws.Range("E32").Value = '  -7.07%  '
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").Value = '4.98'
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("D35").Value = '0.0688'
$ws.Range("E35").Value = '  -1.82%  '
$ws.Range("D36").Value = '4.95'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").Value = '2.54'
$ws.Range("E37").Value = '  +9.49%  '
$ws.Range("D38").Value = '6.47'
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  -5.06%  '
$ws.Range("D40").Value = '0.0251'
$ws.Range("E40").Value = '  -2.69%  '
$ws.Range("E41").Value = '  +8.99%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '18.05'
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("E44").Value = '  +7.98%  '
$ws.Range("D45").Value = '99.08'
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").Value = '0.0946'
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").Value = '1.436.82'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").Value = '0.000204'
$ws.Range("E50").Value = '  -11.33%  '
$ws.Range("D51").Value = '9.82'
$ws.Range("E51").Value = '  -6.02%  '
